{"js": "// Report template refresh: retitle the header, restamp the date, and\n// shorten the footer's confidentiality notice (drop the \"Audix\" brand).\n//\n// The document uses a single section whose \"Primary\" header/footer pair\n// carries the visible text (even-page / first-page variants are untouched\n// by this change, matching the source diff).\nconst section = context.document.sections.getFirst();\nconst header = section.getHeader(\"Primary\");\nconst footer = section.getFooter(\"Primary\");\n\n// --- Header: \"Gain Coin Analysis\" -> \"Audio Bench Report\" -------------\nconst titleHits = header.search(\"Gain Coin Analysis\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"Audio Bench Report\", \"Replace\");\n  await context.sync();\n}\n\n// --- Header: \"February 11, 2025\" -> \"November 21, 2025\" ---------------\nconst dateHits = header.search(\"February 11, 2025\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"November 21, 2025\", \"Replace\");\n  await context.sync();\n}\n\n// --- Footer: \"Audix Confidential\" -> \"Confidential\" --------------------\nconst footerHits = footer.search(\"Audix Confidential\", { matchCase: true });\nfooterHits.load(\"items\");\nawait context.sync();\nif (footerHits.items.length > 0) {\n  footerHits.items[0].insertText(\"Confidential\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Report template refresh: retitle the header, restamp the date, and\n# shorten the footer's confidentiality notice (drop the \"Audix\" brand).\n#\n# The document has a single section; its \"Primary\" (index 1) header/footer\n# pair carries the visible text (even-page / first-page variants are\n# untouched by this change, matching the source diff).\n$d = $word.ActiveDocument\n$sec = $d.Sections(1)\n\n# wdHeaderFooterPrimary = 1\n$hdr = $sec.Headers(1)\n$ftr = $sec.Footers(1)\n\n# --- Header: \"Gain Coin Analysis\" -> \"Audio Bench Report\" -------------\n$hdr.Range.Find.Execute(\n    \"Gain Coin Analysis\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Audio Bench Report\",\n    2\n)\n\n# --- Header: \"February 11, 2025\" -> \"November 21, 2025\" ---------------\n$hdr.Range.Find.Execute(\n    \"February 11, 2025\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"November 21, 2025\",\n    2\n)\n\n# --- Footer: \"Audix Confidential\" -> \"Confidential\" ---------------------\n$ftr.Range.Find.Execute(\n    \"Audix Confidential\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Confidential\",\n    2\n)\n"}
